$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "92.763.86"
$ws.Range("E2").Value = "  -1.64%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.381.89"
$ws.Range("E3").Value = "  -1.24%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.05"
$ws.Range("E5").Value = "  -2.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.54"
$ws.Range("E6").Value = "  -3.52%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.37"
$ws.Range("E7").Value = "  -4.81%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.390"
$ws.Range("E8").Value = "  -3.51%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.963"
$ws.Range("E10").Value = "  -0.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.379.86"
$ws.Range("E11").Value = "  -1.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.94"
$ws.Range("E12").Value = "  +3.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.197"
$ws.Range("E13").Value = "  -0.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.24"
$ws.Range("E14").Value = "  +1.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.606.05"
$ws.Range("E15").Value = "  -1.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.013.91"
$ws.Range("E16").Value = "  -1.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000245"
$ws.Range("E17").Value = "  -2.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.11"
$ws.Range("E18").Value = "  -1.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.374.60"
$ws.Range("E19").Value = "  -1.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.63"
$ws.Range("E20").Value = "  +1.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.46"
$ws.Range("E21").Value = "  -0.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.35"
$ws.Range("E22").Value = "  +3.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "496.16"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.434"
$ws.Range("E24").Value = "  -12.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.53"
$ws.Range("E25").Value = "  +1.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000183"
$ws.Range("E26").Value = "  -4.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "90.68"
$ws.Range("E27").Value = "  -3.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.98"
$ws.Range("E28").Value = "  -0.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.552.16"
$ws.Range("E29").Value = "  -1.53%  "

$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.16"
$ws.Range("E31").Value = "  -3.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.136"
$ws.Range("E32").Value = "  -1.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.68"
$ws.Range("E33").Value = "  -2.47%  "

$ws.Range("E34").Value = "  +0.22%  "

$ws.Range("E35").Value = "  -3.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "29.07"
$ws.Range("E36").Value = "  -2.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.532"
$ws.Range("E37").Value = "  -3.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "557.72"
$ws.Range("E38").Value = "  +2.64%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.47"
$ws.Range("E39").Value = "  -1.96%  "

$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("E41").Value = "  -0.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.38"
$ws.Range("E42").Value = "  -4.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.901"
$ws.Range("E43").Value = "  -0.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.66"
$ws.Range("E44").Value = "  -1.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.71"
$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.66"
$ws.Range("E46").Value = "  +2.80%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0406"
$ws.Range("E47").Value = "  -0.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.43"
$ws.Range("E48").Value = "  -1.96%  "

$ws.Range("B49").Value = "Fantom"
$ws.Range("C49").Value = "https://coinranking.com/coin/uIEWfMFnQo9K_+fantom-ftm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.13"
$ws.Range("E49").Value = "  +19.08%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.32"
$ws.Range("E50").Value = "  -3.10%  "

$ws.Range("E51").Value = "  -3.43%  "
